$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the two consecutive empty "Body Text" paragraphs that sit
#    between the "自动摘要" paragraph and the "Usage" Heading2.
# ---------------------------------------------------------------------
for ($i = 2; $i -lt $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Body Text" -and $p.Range.Text.Trim() -eq "") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Style.NameLocal -eq "Body Text" -and $next.Range.Text.Trim() -eq "") {
            $prev = $d.Paragraphs.Item($i - 1)
            if ($prev.Range.Text.Trim() -eq "自动摘要") {
                $next.Range.Delete()
                $p.Range.Delete()
                break
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Insert discourse connectives ("虽然，" / "然而，") into the sample
#    JSON "content" field to exercise noise filtering.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("为啥差距这么大？去年至今，", $true, $false, $false, $false, $false, `
    $true, 1, $false, "为啥差距这么大？虽然，去年至今，", 2)

$d.Content.Find.Execute("占总贷款余额的54%。但不同地区", $true, $false, $false, $false, $false, `
    $true, 1, $false, "占总贷款余额的54%。然而，但不同地区", 2)

# ---------------------------------------------------------------------
# 3) Re-indent the closing brace of the "request" JSON body sample
#    (the one that follows `"rate": 140`) so it reads "    }".
#    (That paragraph is the only "Source Code" paragraph holding the
#    `"rate": 140` line; the line break before the brace means we can't
#    match the two across a single wildcard, so key off the style +
#    the preceding line instead, then scope Find to that paragraph.)
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Source Code" -and $p.Range.Text -like '*rate*: 140*') {
        $r = $p.Range
        $r.Find.Execute("}", $true, $false, $false, $false, $false, `
            $true, 1, $false, "    }", 2)
        break
    }
}
